# Commit: "Working con proyecto Algerri"
# Rename the two "Vol N Algerri" sheets to "Config N Algerri" and make the
# second config sheet ("Config 2 Algerri") the active sheet/tab, moving the
# selection on that sheet to C21.

$wb = $excel.ActiveWorkbook

# Rename sheets 2 and 3 (by position) from "Vol N Algerri" to "Config N Algerri"
$wb.Worksheets.Item(2).Name = "Config 1 Algerri"
$wb.Worksheets.Item(3).Name = "Config 2 Algerri"

# Make "Config 2 Algerri" the active sheet (this also moves tabSelected from
# the previously-active sheet to this one) and move its selection to C21.
$ws3 = $wb.Worksheets.Item(3)
$ws3.Activate() | Out-Null
$ws3.Range("C21").Select() | Out-Null
